# Apply the data-refresh edit: H2 likes bump + row rotations for the
# feature-extraction dataset (rows 5/6/8 rotate as a 3-cycle, rows
# 13-18 rotate as a 6-cycle). Columns A,E,J,K,L,P,Q are untouched;
# only B,C,D,F,G,H,I,M,N,O move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- H2: likes 309 -> 310 ---
$ws.Range("H2").Value = 310

# --- Row 5 (was Alibaba-NLP/gte-large-en-v1.5) becomes WhereIsAI/UAE-Large-V1 ---
$ws.Range("B5").Value = "WhereIsAI"
$ws.Range("C5").Value = "{'avatarUrl': 'https://www.gravatar.com/avatar/e81bd32cb5ee88835824ad6b60d05697?d=retro&size=100', 'fullname': 'WhereIsAI', 'name': 'WhereIsAI', 'type': 'org', 'isHf': False, 'isEnterprise': False}"
$ws.Range("D5").Value = 277613
$ws.Range("F5").Value = "WhereIsAI/UAE-Large-V1"
$ws.Range("G5").Value = "2024-05-03T02:31:54.000Z"
$ws.Range("H5").Value = 177
$ws.Range("I5").Value = "feature-extraction"
$ws.Range("M5").Value = "1.34GB | 669MB | 337MB"
$ws.Range("N5").Value = 353370112
$ws.Range("O5").Value = 1438814044.16

# --- Row 6 (was WhereIsAI/UAE-Large-V1) becomes Snowflake/snowflake-arctic-embed-l ---
$ws.Range("B6").Value = "Snowflake"
$ws.Range("C6").Value = "{'avatarUrl': 'https://cdn-avatars.huggingface.co/v1/production/uploads/64ba2f59a6ccf0f64b4ad254/eTDA37yFwUVP45c1WTSs2.png', 'fullname': 'Snowflake', 'name': 'Snowflake', 'type': 'org', 'isHf': False, 'isEnterprise': False}"
$ws.Range("D6").Value = 33257
$ws.Range("F6").Value = "Snowflake/snowflake-arctic-embed-l"
$ws.Range("G6").Value = "2024-04-18T19:58:11.000Z"
$ws.Range("H6").Value = 58
$ws.Range("I6").Value = "sentence-similarity"
$ws.Range("M6").Value = "1.34GB | 299MB | 669MB | 337MB | 318MB | 337MB | 337MB"
$ws.Range("N6").Value = 313524224
# O6 stays 1438814044.16 (unchanged by the diff)

# --- Row 8 (was Snowflake/snowflake-arctic-embed-l) becomes Alibaba-NLP/gte-large-en-v1.5 ---
$ws.Range("B8").Value = "Alibaba-NLP"
$ws.Range("C8").Value = "{'avatarUrl': 'https://www.gravatar.com/avatar/1ae3fd9f5b9356f196c997d93eb23038?d=retro&size=100', 'fullname': 'Alibaba-NLP', 'name': 'Alibaba-NLP', 'type': 'org', 'isHf': False, 'isEnterprise': False}"
$ws.Range("D8").Value = 75109
$ws.Range("F8").Value = "Alibaba-NLP/gte-large-en-v1.5"
$ws.Range("G8").Value = "2024-04-26T13:51:26.000Z"
$ws.Range("H8").Value = 56
# I8 stays sentence-similarity (unchanged by the diff)
$ws.Range("M8").Value = "1.75GB | 361MB | 873MB | 446MB | 387MB | 446MB | 446MB"
$ws.Range("N8").Value = 378535936
$ws.Range("O8").Value = 1879048192

# --- Row 13 (was Alibaba-NLP/gte-base-en-v1.5) becomes Xenova/all-MiniLM-L6-v2 ---
$ws.Range("B13").Value = "Xenova"
$ws.Range("C13").Value = "{'avatarUrl': 'https://cdn-avatars.huggingface.co/v1/production/uploads/61b253b7ac5ecaae3d1efe0c/hwiQ0uvz3t-L5a-NtBIO6.png', 'fullname': 'Joshua', 'name': 'Xenova', 'type': 'user', 'isPro': False, 'isHf': True}"
$ws.Range("D13").Value = 105
$ws.Range("F13").Value = "Xenova/all-MiniLM-L6-v2"
$ws.Range("G13").Value = "2024-03-12T03:11:08.000Z"
$ws.Range("H13").Value = 39
$ws.Range("I13").Value = "feature-extraction"
$ws.Range("M13").Value = "90.4MB | 45.3MB | 23MB | 23MB"
$ws.Range("N13").Value = 24117248
$ws.Range("O13").Value = 94791270.40000001

# --- Row 14 (was Xenova/all-MiniLM-L6-v2) becomes Xenova/larger_clap_music_and_speech ---
# B14/C14 stay Xenova / Joshua avatar (unchanged by the diff)
$ws.Range("D14").Value = 3
$ws.Range("F14").Value = "Xenova/larger_clap_music_and_speech"
$ws.Range("G14").Value = "2024-03-17T15:10:25.000Z"
$ws.Range("H14").Value = 2
# I14 stays feature-extraction (unchanged by the diff)
$ws.Range("M14").Value = "783MB | 395MB | 205MB"
$ws.Range("N14").Value = 214958080
$ws.Range("O14").Value = 821035008

# --- Row 15 (was Xenova/larger_clap_music_and_speech) becomes jinaai/jina-embeddings-v2-base-de ---
$ws.Range("B15").Value = "jinaai"
$ws.Range("C15").Value = "{'avatarUrl': 'https://cdn-avatars.huggingface.co/v1/production/uploads/603763514de52ff951d89793/AFoybzd5lpBQXEBrQHuTt.png', 'fullname': 'Jina AI', 'name': 'jinaai', 'type': 'org', 'isHf': False, 'isEnterprise': True}"
$ws.Range("D15").Value = 20417
$ws.Range("F15").Value = "jinaai/jina-embeddings-v2-base-de"
$ws.Range("G15").Value = "2024-04-26T07:35:44.000Z"
$ws.Range("H15").Value = 50
# I15 stays feature-extraction (unchanged by the diff)
$ws.Range("M15").Value = "641MB | 321MB | 162MB"
$ws.Range("N15").Value = 169869312
$ws.Range("O15").Value = 672137216

# --- Row 16 (was jinaai/jina-embeddings-v2-base-de) becomes mixedbread-ai/mxbai-embed-2d-large-v1 ---
$ws.Range("B16").Value = "mixedbread-ai"
$ws.Range("C16").Value = "{'avatarUrl': 'https://cdn-avatars.huggingface.co/v1/production/uploads/643ee0870d1194da249bd7fe/voYdYlFgQH5vyMwyMNluZ.png', 'fullname': 'mixedbread ai', 'name': 'mixedbread-ai', 'type': 'org', 'isHf': False, 'isEnterprise': False}"
$ws.Range("D16").Value = 35205
$ws.Range("F16").Value = "mixedbread-ai/mxbai-embed-2d-large-v1"
$ws.Range("G16").Value = "2024-04-04T21:36:56.000Z"
$ws.Range("H16").Value = 29
# I16 stays feature-extraction (unchanged by the diff)
$ws.Range("M16").Value = "1.34GB | 669MB | 337MB"
$ws.Range("N16").Value = 353370112
$ws.Range("O16").Value = 1438814044.16

# --- Row 17 (was mixedbread-ai/mxbai-embed-2d-large-v1) becomes Xenova/text2vec-base-chinese-paraphrase ---
$ws.Range("B17").Value = "Xenova"
$ws.Range("C17").Value = "{'avatarUrl': 'https://cdn-avatars.huggingface.co/v1/production/uploads/61b253b7ac5ecaae3d1efe0c/hwiQ0uvz3t-L5a-NtBIO6.png', 'fullname': 'Joshua', 'name': 'Xenova', 'type': 'user', 'isPro': False, 'isHf': True}"
$ws.Range("D17").Value = 3
$ws.Range("F17").Value = "Xenova/text2vec-base-chinese-paraphrase"
$ws.Range("G17").Value = "2024-03-23T00:34:15.000Z"
$ws.Range("H17").Value = 1
# I17 stays feature-extraction (unchanged by the diff)
$ws.Range("M17").Value = "470MB | 178MB | 235MB | 119MB | 183MB | 119MB | 119MB"
$ws.Range("N17").Value = 124780544
$ws.Range("O17").Value = 492830720

# --- Row 18 (was Xenova/text2vec-base-chinese-paraphrase) becomes Alibaba-NLP/gte-base-en-v1.5 ---
$ws.Range("B18").Value = "Alibaba-NLP"
$ws.Range("C18").Value = "{'avatarUrl': 'https://www.gravatar.com/avatar/1ae3fd9f5b9356f196c997d93eb23038?d=retro&size=100', 'fullname': 'Alibaba-NLP', 'name': 'Alibaba-NLP', 'type': 'org', 'isHf': False, 'isEnterprise': False}"
$ws.Range("D18").Value = 97301
$ws.Range("F18").Value = "Alibaba-NLP/gte-base-en-v1.5"
$ws.Range("G18").Value = "2024-04-26T13:53:41.000Z"
$ws.Range("H18").Value = 9
$ws.Range("I18").Value = "sentence-similarity"
$ws.Range("M18").Value = "556MB | 167MB | 278MB | 147MB | 174MB | 147MB | 147MB"
$ws.Range("N18").Value = 154140672
$ws.Range("O18").Value = 583008256
